# Update NATMI TPM recomputation results (new ligand-expressing cell count)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.319712
$ws.Range("H2").Value = 3.959136
$ws.Range("M2").Value = 3.046487666666666
$ws.Range("N2").Value = 9.139462999999999
$ws.Range("O2").Value = 0.01743556347089316
$ws.Range("P2").Value = 0.01743556347089316
$ws.Range("Q2").Value = 4.020486331551999
$ws.Range("R2").Value = 36.18437698396799
$ws.Range("S2").Value = 0.01743556347089316
$ws.Range("T2").Value = 0.01743556347089316

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.319712
$ws.Range("H3").Value = 3.959136
$ws.Range("M3").Value = 154.6652373333333
$ws.Range("O3").Value = 0.8851752763590445
$ws.Range("P3").Value = 0.8851752763590445
$ws.Range("Q3").Value = 204.113569691648
$ws.Range("R3").Value = 1837.022127224832
$ws.Range("S3").Value = 0.8851752763590445
$ws.Range("T3").Value = 0.8851752763590445

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.319712
$ws.Range("H4").Value = 3.959136
$ws.Range("M4").Value = 15.82721166666667
$ws.Range("N4").Value = 47.481635
$ws.Range("O4").Value = 0.09058180559889371
$ws.Range("P4").Value = 0.09058180559889373
$ws.Range("Q4").Value = 20.88736116304
$ws.Range("R4").Value = 187.98625046736
$ws.Range("S4").Value = 0.09058180559889371
$ws.Range("T4").Value = 0.09058180559889373

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.319712
$ws.Range("H5").Value = 3.959136
$ws.Range("M5").Value = 1.189438
$ws.Range("N5").Value = 3.568314
$ws.Range("O5").Value = 0.006807354571168639
$ws.Range("P5").Value = 0.00680735457116864
$ws.Range("Q5").Value = 1.569715601856
$ws.Range("R5").Value = 14.127440416704
$ws.Range("S5").Value = 0.006807354571168639
$ws.Range("T5").Value = 0.00680735457116864
